# Localization status report regeneration: flip the "Ready for handoff"
# status to "In Translation" everywhere it appears, and re-fit the
# status columns now that the text is shorter.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# New status text is shorter than the old one, so the status columns can
# shrink to fit. ColumnWidth is in character units on a whole-pixel grid
# (pixels = round(width*6), stored width = (pixels+5)/6); 12.5 is the
# closest achievable setting to the desired ~13.41 character width.
$newColumnWidth = 12.5

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange

    foreach ($cell in $used.Cells) {
        # NB: keep the literal on the LEFT of -eq. Cells that hold the text
        # "True"/"False" come back from Value() as real [bool]s, and with a
        # [bool] on the left PowerShell coerces the right-hand string to a
        # bool too (any non-empty string -> $true), producing false
        # positives. A string literal on the left forces a string compare.
        if ($oldStatus -eq $cell.Value()) {
            $cell.Value = $newStatus
        }
    }

    # Resize every column whose header is "Status" (zh-cn/de-de detail
    # sheets) or that holds the status text directly (Overview's zh-cn /
    # de-de columns).
    $colCount = $used.Columns.Count
    for ($c = 1; $c -le $colCount; $c++) {
        $header = $ws.Cells.Item(1, $c).Value()
        if ("Status" -eq $header -or "zh-cn" -eq $header -or "de-de" -eq $header) {
            $ws.Columns.Item($c).ColumnWidth = $newColumnWidth
        }
    }
}
